# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (currently empty) column
# inserted just before the old "Late" column, pushing the old
# N/O/P ("Late" / "heading" / "Outstanding") columns one place to the
# right (-> O/P/Q). The new column inherits the width of the column to
# its left, same as Excel does when you insert a column from the UI.
#
# The workbook's active tab also moves from "Summary" to
# "Repayment schedule", with the selected cell on that sheet moving to R6.

$wb = $excel.ActiveWorkbook

$repaymentSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new, blank column before column N (14th column == "Late").
# This shifts the existing N/O/P columns to O/P/Q and grows the sheet's
# dimension from A1:P16 to A1:Q16 automatically.
$repaymentSchedule.Columns.Item(14).Insert()

# Excel gives a freshly inserted column the width of its left neighbour
# (column M here); set it explicitly so the stored width matches (~11
# characters, same as column M).
$repaymentSchedule.Columns.Item(14).ColumnWidth = 10.1

# Make "Repayment schedule" the active/selected sheet, with R6 selected,
# which also clears the previous tab selection on "Summary".
$repaymentSchedule.Activate()
$repaymentSchedule.Range("R6").Select()
